# Update the "expected_value" cell for the first test case row from
# " insulin, tylenol, vaccine" to " APIXABAN, ETANERCEPT".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = " APIXABAN, ETANERCEPT"
